$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

$ws.Range("D2").Value = "26.167.88"
$ws.Range("E2").Value = "  -1.32%  "
$ws.Range("D3").Value = "1.653.92"
$ws.Range("E3").Value = "  -1.70%  "
Set-TextCell "D5" "218.44"
$ws.Range("E5").Value = "  +0.34%  "
Set-TextCell "D6" "0.5211"
$ws.Range("E6").Value = "  -2.16%  "
$ws.Range("E7").Value = "  +0.30%  "
Set-TextCell "D8" "0.2668"
$ws.Range("E8").Value = "  -0.34%  "
Set-TextCell "D9" "0.06327"
$ws.Range("E9").Value = "  -1.57%  "
Set-TextCell "D10" "21.13"
$ws.Range("E10").Value = "  -1.69%  "
Set-TextCell "D11" "0.07716"
$ws.Range("E11").Value = "  -1.15%  "
$ws.Range("D12").Value = "1.655.97"
$ws.Range("E12").Value = "  -1.62%  "
Set-TextCell "D13" "4.438"
$ws.Range("E13").Value = "  -1.65%  "
$ws.Range("D14").Value = "1.880.31"
$ws.Range("E14").Value = "  -1.68%  "
Set-TextCell "D15" "0.5472"
$ws.Range("E15").Value = "  -2.72%  "
$ws.Range("D16").Value = "0.0₅8254"
$ws.Range("E16").Value = "  -2.06%  "
Set-TextCell "D17" "64.89"
$ws.Range("E17").Value = "  -1.71%  "
$ws.Range("D18").Value = "26.173.77"
$ws.Range("E18").Value = "  -1.41%  "
$ws.Range("E19").Value = "  +0.36%  "
Set-TextCell "D20" "4.671"
$ws.Range("E20").Value = "  -2.85%  "
Set-TextCell "D21" "193.06"
$ws.Range("E21").Value = "  -1.42%  "
Set-TextCell "D22" "10.17"
$ws.Range("E22").Value = "  -2.30%  "
Set-TextCell "D23" "6.101"
$ws.Range("E23").Value = "  -4.47%  "
$ws.Range("E24").Value = "  +0.47%  "
Set-TextCell "D25" "137.72"
$ws.Range("E25").Value = "  -3.85%  "
$ws.Range("E26").Value = "  -2.83%  "
Set-TextCell "D27" "7.239"
$ws.Range("E27").Value = "  -3.13%  "
Set-TextCell "D28" "16.15"
$ws.Range("E28").Value = "  -0.22%  "
Set-TextCell "D29" "1.431"
$ws.Range("E29").Value = "  +1.02%  "
Set-TextCell "D30" "0.06036"
Set-TextCell "D31" "1.283"
$ws.Range("E31").Value = "  +0.28%  "
Set-TextCell "D32" "3.564"
$ws.Range("E32").Value = "  -1.18%  "
Set-TextCell "D33" "3.336"
$ws.Range("E33").Value = "  -3.54%  "
Set-TextCell "D34" "1.650"
$ws.Range("E34").Value = "  -3.22%  "
Set-TextCell "D35" "0.9817"
$ws.Range("E35").Value = "  -3.40%  "
$ws.Range("E36").Value = "  -0.37%  "
Set-TextCell "D37" "2.767"
$ws.Range("E37").Value = "  -0.85%  "
Set-TextCell "D38" "0.5928"
$ws.Range("E38").Value = "  +3.99%  "
Set-TextCell "D39" "0.01592"
$ws.Range("E39").Value = "  -2.95%  "
Set-TextCell "D40" "5.968"
$ws.Range("E40").Value = "  +0.19%  "
Set-TextCell "D41" "0.8630"
$ws.Range("E41").Value = "  -0.86%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.044.09"
$ws.Range("E42").Value = "  -1.75%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell "D43" "1.003"
$ws.Range("E43").Value = "  +0.16%  "
Set-TextCell "D44" "99.62"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("D45").Value = "1.793.72"
$ws.Range("E45").Value = "  -2.13%  "
$ws.Range("D46").Value = "0.0₈113"
$ws.Range("E46").Value = "  +1.29%  "
Set-TextCell "D47" "57.24"
$ws.Range("E47").Value = "  -0.07%  "
Set-TextCell "D48" "1.006"
$ws.Range("E48").Value = "  +0.72%  "
Set-TextCell "D49" "8.108"
$ws.Range("E49").Value = "  -0.45%  "
Set-TextCell "D50" "0.05178"
$ws.Range("E50").Value = "  -0.48%  "
Set-TextCell "D51" "1.467"
$ws.Range("E51").Value = "  +3.85%  "
